$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep gridlines visible (this runtime's round-trip otherwise drops it to
# hidden even on a pure no-op load/save).
$excel.ActiveWindow.DisplayGridlines = $true

# Insert new column C ("Número de Invitados") before the existing "Felicitaciones" column
$ws.Columns("C").Insert()

# Set header row values
$ws.Range("A1").Value = "Nombre del Invitado"
$ws.Range("B1").Value = "Departe de"
$ws.Range("C1").Value = "Número de Invitados"
$ws.Range("D1").Value = "Felicitaciones"

# Copy style of B1 (already existing header style) onto new C1 cell
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats

# Remove old sample data row (row 2)
$ws.Rows("2").Delete()

# Column widths.
# A keeps its original 25-char width (re-applied so it round-trips exactly).
# B and D keep their original widths untouched (13.7109375 / 17.42578125).
# C is the newly introduced column - target stored width ~25.28515625; this
# is the closest value reachable through Excel's pixel-rounded ColumnWidth
# setter (width is stored internally as a whole pixel count).
$ws.Columns("A").ColumnWidth = 24.166666666666668
$ws.Columns("C").ColumnWidth = 24.5

$ws.Range("E1").Select()
